$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 276; this shifts the existing rows 276-341 down
# to 277-342 (so the old row 341 becomes the new row 342), matching the
# target diff exactly.
$ws.Rows.Item(276).Insert()

# Populate the newly inserted row 276 with the new record.
$ws.Cells.Item(276, 1).Value = 5
$ws.Cells.Item(276, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(276, 3).Value = "Maule"
$ws.Cells.Item(276, 4).Value = 44722
$ws.Cells.Item(276, 5).Value = 7
$ws.Cells.Item(276, 6).Value = 100112023
$ws.Cells.Item(276, 7).Value = "Brócoli"
$ws.Cells.Item(276, 8).Value = "Sin especificar"
$ws.Cells.Item(276, 9).Value = "Primera"
$ws.Cells.Item(276, 10).Value = 3000
$ws.Cells.Item(276, 11).Value = 1000
$ws.Cells.Item(276, 12).Value = 1000
$ws.Cells.Item(276, 13).Value = 1000
$ws.Cells.Item(276, 14).Value = "$/unidad"
$ws.Cells.Item(276, 15).Value = "Región del Maule"
$ws.Cells.Item(276, 16).Value = 1000
$ws.Cells.Item(276, 17).Value = 1
$ws.Cells.Item(276, 18).Value = "Hortaliza"

# Match the date-formatted style used by the other rows in column D.
$ws.Cells.Item(276, 4).NumberFormat = $ws.Cells.Item(277, 4).NumberFormat
